# Timesheet update — log additional time for the homing/target-assist
# missile work and the "enemies can be damaged" feature on 45310 (row 27):
#   - LO1 (utilities, col B): +6 more minutes on top of the existing 6
#   - LO2 (Character + NPC, col C): a new 5-minute entry (was untouched/0)
#   - LO4 (world interactions, col E): three more entries (26+26+15 min)
# Column D (LO3) and every other row/date are left untouched. Row totals
# (F), grand totals (row 3) and percentages (row 4) all recalc from these.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B27").Formula = "=(1/60)*(6+6)"
$ws.Range("C27").Formula = "=(1/60)*(5)"
$ws.Range("E27").Formula = "=(1/60)*(21+26+20+20+26+26+15)"

# Leave the selection on D27, matching where the edit was being made.
$ws.Range("D27").Select()
